# Báo cáo (BCT3_Dung.docx): sửa ngày tháng của kỳ đánh giá.
# Kỳ đánh giá đổi từ "6/09/2011 đến 13/09/2011" thành "6/10/2011 đến 13/10/2011"
# (tháng 09 -> tháng 10), cả hai lần xuất hiện của "/09" trong văn bản.

$d = $word.ActiveDocument

# Thay toàn bộ "/09" bằng "/10" (chỉ xuất hiện trong câu "Kỳđánhgiá",
# hai lần - ngày bắt đầu và ngày kết thúc).
$d.Content.Find.Execute("/09", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "/10", 2) | Out-Null
